# plotEIC methods for fGroupsSet
# - add a new "getEICsForFGroups" row above "getFeatures" (row 19), marked as
#   implemented (D), using ionization (F) and done (G)
# - mark plotEIC (row 34 after the insert) as done (G) in addition to implement (D)
# - move the active-cell selection down to G35 to reflect the new cursor position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 19, pushing getFeatures..groupAlgorithm down by one
$ws.Rows.Item(19).Insert() | Out-Null

# Populate the new row for the getEICsForFGroups method
$ws.Range("A19").Value = "getEICsForFGroups"
$ws.Range("D19").Value = "X"
$ws.Range("F19").Value = "X"
$ws.Range("G19").Value = "X"

# plotEIC (originally row 33, now row 34) is now also marked as done
$ws.Range("G34").Value = "X"

# Update the selected cell to match the saved view state
$ws.Range("G35").Select() | Out-Null
